$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.036.12"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "3.892.59"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "466.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.737"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000334"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "4.516.83"
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "3.895.56"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("D20").Value = "67.243.68"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "38.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.65%  "
$ws.Range("E27").Value = "  +6.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("E29").Value = "  -3.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "735.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "0.0₃0781"
$ws.Range("E38").Value = "  +13.99%  "
$ws.Range("E39").Value = "  -5.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0476"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.337"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.46%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.140"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +5.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.13%  "
$ws.Range("E47").Value = "  -5.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
